# Weekly update: insert two new price records at the top of the
# "Arándano (blue)" data block (rows 280-281), pushing the existing
# rows down by 2 (old row 280 -> new row 282, ... old row 339 -> new row 341).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 280/281; everything from old row 280 onward
# shifts down to rows 282 onward.
$ws.Range("A280:A281").EntireRow.Insert()

# Seed the two new rows by duplicating the (now shifted) rows that follow
# them, then overwrite just the cells that carry genuinely new data. This
# mirrors the original records' shared columns (market/product metadata)
# without retyping every field.
$ws.Rows.Item(282).Copy()
$ws.Rows.Item(280).PasteSpecial()

$ws.Rows.Item(283).Copy()
$ws.Rows.Item(281).PasteSpecial()

$excel.CutCopyMode = 0

# New row 280 (Provincia de Curicó, Primera) — price jump for the latest week.
$ws.Cells.Item(280, 4).Value2 = 45218
$ws.Cells.Item(280, 14).Value2 = 10000
$ws.Cells.Item(280, 15).Value2 = 10000
$ws.Cells.Item(280, 16).Value2 = 10000
$ws.Cells.Item(280, 19).Value2 = 5000

# New row 281 (Provincia de Linares, Primera) — new weekly record.
$ws.Cells.Item(281, 4).Value2 = 45218
$ws.Cells.Item(281, 13).Value2 = 160
$ws.Cells.Item(281, 14).Value2 = 10000
$ws.Cells.Item(281, 15).Value2 = 11000
$ws.Cells.Item(281, 16).Value2 = 10500
$ws.Cells.Item(281, 19).Value2 = 5250
